$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.583.68"
$ws.Range("E2").Value = "  +0.13%  "

$ws.Range("D3").Value = "2.470.76"
$ws.Range("E3").Value = "  -0.36%  "

$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.40"
$ws.Range("E5").Value = "  +1.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.94"
$ws.Range("E6").Value = "  -1.04%  "

$ws.Range("E7").Value = "  +1.00%  "

$ws.Range("E8").Value = "  +0.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.514"
$ws.Range("E9").Value = "  +1.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0860"
$ws.Range("E10").Value = "  +9.36%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "32.85"
$ws.Range("E11").Value = "  +0.22%  "

$ws.Range("E12").Value = "  +0.22%  "

$ws.Range("D13").Value = "2.854.20"
$ws.Range("E13").Value = "  -0.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.87"
$ws.Range("E14").Value = "  -0.10%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.53"
$ws.Range("E15").Value = "  -4.45%  "

$ws.Range("D16").Value = "2.451.48"
$ws.Range("E16").Value = "  -2.37%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.785"
$ws.Range("E17").Value = "  +1.96%  "

$ws.Range("D18").Value = "41.527.83"
$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0949"
$ws.Range("E19").Value = "  -0.01%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.45"
$ws.Range("E20").Value = "  -1.79%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.16"
$ws.Range("E21").Value = "  -1.63%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.30"
$ws.Range("E22").Value = "  +0.38%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.78"
$ws.Range("E23").Value = "  +0.96%  "

$ws.Range("E24").Value = "  +0.84%  "

$ws.Range("E25").Value = "  +1.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.60"
$ws.Range("E27").Value = "  -0.91%  "

$ws.Range("E28").Value = "  +2.96%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.85"
$ws.Range("E29").Value = "  +1.74%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.19"
$ws.Range("E30").Value = "  +0.92%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.69"
$ws.Range("E31").Value = "  +1.59%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.49"
$ws.Range("E32").Value = "  +0.79%  "

$ws.Range("E33").Value = "  +0.12%  "

$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.58"
$ws.Range("E34").Value = "  +0.39%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0766"
$ws.Range("E35").Value = "  +1.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.15"
$ws.Range("E36").Value = "  -1.95%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.89"
$ws.Range("E37").Value = "  -0.59%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.115"
$ws.Range("E38").Value = "  +1.18%  "

$ws.Range("E39").Value = "  +0.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.102"
$ws.Range("E40").Value = "  -3.49%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.95"
$ws.Range("E41").Value = "  -3.46%  "

$ws.Range("E42").Value = "  +2.90%  "

$ws.Range("D43").Value = "1.987.59"
$ws.Range("E43").Value = "  +1.30%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0284"
$ws.Range("E44").Value = "  +0.19%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.98"
$ws.Range("E45").Value = "  -1.92%  "

$ws.Range("E46").Value = "  +0.77%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.18"
$ws.Range("E47").Value = "  +2.08%  "

$ws.Range("D48").Value = "2.711.51"
$ws.Range("E48").Value = "  -0.22%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.47"
$ws.Range("E49").Value = "  -0.18%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.62"
$ws.Range("E50").Value = "  +1.43%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "66.65"
$ws.Range("E51").Value = "  -2.24%  "
